$d = $word.ActiveDocument

# 1. Insert the new clause ", the Fast ESP Wi-Fi sample," after "ESP-IDF" in the
#    paragraph describing creating a new project.
$d.Content.Find.Execute(
    "Create a new project using ESP-IDF and make sure",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Create a new project using ESP-IDF, the Fast ESP Wi-Fi sample, and make sure",
    2
)

# 2. Explicitly set the page orientation to portrait (adds w:orient="portrait"
#    to the section's pgSz element).
$sec = $d.Sections.Item(1)
$sec.PageSetup.Orientation = 0
